$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting
# (some updated values would otherwise be auto-converted to numbers by Excel)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "63.033.08"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "3.270.78"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "599.07"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").Value = "138.28"
$ws.Range("E6").Value = "  -2.36%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.269.37"
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("E10").Value = "  -0.42%  "
$ws.Range("D11").Value = "5.45"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").Value = "0.461"
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("D13").Value = "0.0000242"
$ws.Range("E13").Value = "  -2.45%  "
$ws.Range("D14").Value = "33.97"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").Value = "3.807.57"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("D17").Value = "3.273.03"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").Value = "63.123.48"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").Value = "6.74"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").Value = "471.98"
$ws.Range("E20").Value = "  -1.37%  "
$ws.Range("D21").Value = "13.80"
$ws.Range("E21").Value = "  -3.07%  "
$ws.Range("D22").Value = "0.725"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("D23").Value = "7.85"
$ws.Range("E23").Value = "  -1.78%  "
$ws.Range("D24").Value = "13.64"
$ws.Range("E24").Value = "  +2.59%  "
$ws.Range("D25").Value = "84.09"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").Value = "7.07"
$ws.Range("E29").Value = "  -2.18%  "
$ws.Range("D30").Value = "7.95"
$ws.Range("E30").Value = "  -1.98%  "
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("D32").Value = "28.16"
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("D33").Value = "0.104"
$ws.Range("E33").Value = "  -2.96%  "
$ws.Range("E34").Value = "  -3.75%  "
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("D36").Value = "5.92"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").Value = "51.76"
$ws.Range("E37").Value = "  -1.85%  "
$ws.Range("D38").Value = "0.0₃0719"
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("D40").Value = "3.084.29"
$ws.Range("E40").Value = "  +2.69%  "
$ws.Range("D41").Value = "421.70"
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("D42").Value = "0.117"
$ws.Range("E42").Value = "  +6.10%  "
$ws.Range("D43").Value = "8.21"
$ws.Range("E43").Value = "  -2.19%  "
$ws.Range("E44").Value = "  -4.49%  "
$ws.Range("E45").Value = "  -3.30%  "
$ws.Range("D46").Value = "2.17"
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("D48").Value = "126.80"
$ws.Range("E48").Value = "  +3.17%  "
$ws.Range("D49").Value = "35.67"
$ws.Range("E49").Value = "  +5.53%  "
$ws.Range("D50").Value = "25.82"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("E51").Value = "  -1.81%  "
